$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$c = $t.Cell(1, 1)
$c.Range.Text = "134÷8="
$c = $t.Cell(1, 2)
$c.Range.Text = "186÷5="
$c = $t.Cell(1, 3)
$c.Range.Text = "653÷3="
$c = $t.Cell(1, 4)
$c.Range.Text = "701÷5="
$c = $t.Cell(1, 5)
$c.Range.Text = "771÷4="
$c = $t.Cell(5, 1)
$c.Range.Text = "379÷3="
$c = $t.Cell(5, 2)
$c.Range.Text = "115÷4="
$c = $t.Cell(5, 3)
$c.Range.Text = "718÷9="
$c = $t.Cell(5, 4)
$c.Range.Text = "407÷7="
$c = $t.Cell(5, 5)
$c.Range.Text = "126÷9="
$c = $t.Cell(9, 1)
$c.Range.Text = "551÷8="
$c = $t.Cell(9, 2)
$c.Range.Text = "608÷5="
$c = $t.Cell(9, 3)
$c.Range.Text = "628÷7="
$c = $t.Cell(9, 4)
$c.Range.Text = "442÷6="
$c = $t.Cell(9, 5)
$c.Range.Text = "268÷8="
$c = $t.Cell(13, 1)
$c.Range.Text = "136÷7="
$c = $t.Cell(13, 2)
$c.Range.Text = "806÷5="
$c = $t.Cell(13, 3)
$c.Range.Text = "163÷8="
$c = $t.Cell(13, 4)
$c.Range.Text = "838÷2="
$c = $t.Cell(13, 5)
$c.Range.Text = "359÷8="
$c = $t.Cell(17, 1)
$c.Range.Text = "609÷6="
$c = $t.Cell(17, 2)
$c.Range.Text = "507÷2="
$c = $t.Cell(17, 3)
$c.Range.Text = "919÷4="
$c = $t.Cell(17, 4)
$c.Range.Text = "833÷8="
$c = $t.Cell(17, 5)
$c.Range.Text = "729÷5="
